# Kd_FP_sum.xlsx — "added one peptide to FP-sum"
#
# A new peptide row (DFATTV) is appended to the bottom of the sorted FP-sum
# table, just above the trailing thick-bottom-border row. Excel's
# "extend formatting" behaviour then carries the new row's formatting one
# row further down (row 37), which is why that row shows up pre-formatted
# but empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data row (row 35) -------------------------------------------------
# A35 already carries the Courier-New / size-11 / centered look used
# throughout column A, so only the values need to be written.
$ws.Range("A35").Value = "DFATTV"
$ws.Range("B35").ClearFormats()
$ws.Range("C35").ClearFormats()
$ws.Range("B35").Value = 108.76
$ws.Range("C35").Value = 9.62

# --- formatting spills one row further down (row 37), matching Excel's ----
# --- "extend data range formats" behaviour ---------------------------------
$ws.Range("A35").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- selection / view state, matching the post-edit cursor position -------
$ws.Range("A37:C37").Select()
$excel.ActiveWindow.ScrollRow = 28
